$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row2 = @(1.95, 2.18, 3.2, 4.2, 3.4, 4.3, 1.33, 1.04, 3.9, 1.24, 2.16, 1.69, 1.46, 2.48, 1.63, 2.28, 1.31, 1.84, 24, 120, 980, 85, 44, 14, 980, 980, 980, 9.4, 90, 980, 900, 980, 980, 580, 55, 980)
$row3 = @(2.5, 2.96, 2.5, 2.94, 3.55, 4, 1.3, 1.04, 4.4, 1.2, 2.16, 1.6, 1.5, 2.48, 1.53, 2.42, 1.53, 1.52, 23, 16.5, 22, 42, 16.5, 9.800000000000001, 13.5, 28, 23, 14, 16, 34, 42, 28, 65, 580, 18, 18)
$row4 = @(1.36, 1.39, 9.4, 14.5, 5.3, 5.9, 0, 0, 0, 0, 2.68, 1.5, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$row5 = @(7.6, 7.8, 1.58, 1.59, 4.1, 4.2, 1.46, 1.08, 3.45, 1.39, 1.81, 2.2, 1.3, 4.1, 2.26, 1.76, 2.68, 1.14, 12.5, 6.8, 7.8, 14, 21, 9.199999999999999, 10, 18.5, 60, 29, 28, 48, 260, 140, 140, 200, 220, 11)
$row6 = @(2.3, 2.4, 3.2, 3.45, 3.5, 3.7, 1.34, 1.06, 4.2, 1.26, 2.08, 1.76, 1.44, 2.92, 1.63, 2.42, 1.4, 1.71, 34, 15.5, 980, 300, 12.5, 8.6, 26, 90, 38, 12.5, 25, 980, 120, 60, 980, 580, 16.5, 980)
$row7 = @(1.49, 1.51, 6.6, 7.4, 5, 5.2, 1.25, 1.03, 6.2, 1.16, 2.78, 1.5, 1.71, 2.2, 1.67, 2.28, 1.16, 2.88, 75, 120, 220, 210, 13.5, 13, 28, 460, 12.5, 11, 22, 360, 15, 15, 28, 300, 5.4, 200)
$row8 = @(2.22, 2.3, 3.8, 4.1, 3.25, 3.3, 1.4, 1.09, 3.1, 1.41, 1.72, 2.22, 1.26, 4.3, 1.92, 1.96, 1.32, 1.76, 12, 12.5, 90, 200, 8.4, 7.6, 16.5, 240, 13.5, 11, 50, 350, 980, 75, 120, 580, 980, 160)
$row9 = @(2.08, 2.14, 4.1, 4.4, 3.35, 3.45, 1.01, 1.09, 3.3, 1.4, 1.84, 2.16, 1.29, 4, 1.93, 1.98, 1.29, 1.87, 12.5, 14, 80, 290, 8.6, 7.8, 18.5, 65, 13, 11.5, 21, 190, 28, 26, 46, 580, 21, 250)
$row10 = @(2.26, 2.34, 3.75, 4, 3.15, 3.3, 1.5, 1.11, 2.86, 1.48, 1.63, 2.42, 1.23, 4.7, 2, 1.88, 1.33, 1.74, 9.6, 11.5, 65, 1000, 8, 7.4, 30, 250, 13.5, 11.5, 60, 350, 120, 85, 130, 580, 1000, 1000)
$row11 = @(4.8, 5.1, 1.92, 1.98, 3.45, 3.6, 1.01, 1.09, 3.35, 1.38, 1.78, 2.14, 1.29, 4, 1.94, 1.95, 2.02, 1.25, 12.5, 8, 11.5, 90, 15.5, 8, 10.5, 70, 980, 36, 38, 980, 140, 440, 350, 580, 1000, 16.5)
$row12 = @(2.4, 2.52, 3.55, 3.8, 3.05, 3.1, 1.55, 1.12, 2.8, 1.5, 1.6, 2.46, 1.22, 5.1, 2.02, 1.87, 1.35, 1.65, 9.4, 11, 26, 80, 8, 7.2, 16, 65, 14, 12, 22, 80, 40, 34, 60, 180, 34, 75)
$row13 = @(2.92, 3.1, 2.64, 2.74, 3.3, 3.45, 1.01, 1.09, 3.1, 1.41, 1.73, 1.99, 1.27, 4.3, 1.74, 2, 1.57, 1.48, 11, 10, 16, 42, 11, 7.4, 12.5, 34, 21, 13.5, 20, 55, 340, 110, 60, 580, 42, 34)
$row14 = @(3.05, 3.25, 2.66, 2.78, 3.1, 3.15, 1.51, 1.1, 2.92, 1.46, 1.65, 2.36, 1.24, 4.7, 1.96, 2, 1.56, 1.44, 10.5, 9, 16.5, 980, 10, 7.4, 13, 980, 21, 14.5, 60, 200, 340, 980, 370, 580, 980, 980)
$row15 = @(1.73, 1.85, 4.8, 5.4, 3.8, 4.3, 1.31, 1.05, 4.8, 1.2, 2.26, 1.66, 1.52, 2.64, 1.6, 2.42, 1.22, 2.16, 980, 980, 980, 130, 980, 9.6, 980, 130, 980, 40, 980, 130, 980, 980, 980, 1000, 29, 1000)
$row16 = @(1.5, 1.54, 6.8, 8.199999999999999, 4.7, 5.2, 1.3, 1.05, 4.1, 1.27, 2.06, 1.81, 1.42, 3.05, 2, 1.89, 1.14, 2.84, 18, 32, 65, 250, 8.6, 11, 29, 120, 8.800000000000001, 9.800000000000001, 26, 110, 13.5, 16, 38, 160, 8, 160)
$row17 = @(2.24, 2.36, 3.4, 3.8, 3.3, 3.55, 1.39, 1.09, 3.05, 1.41, 1.62, 1.99, 1.27, 3.25, 1.9, 1.94, 1.37, 1.73, 21, 22, 980, 300, 8.800000000000001, 12.5, 18, 980, 32, 18, 60, 170, 980, 980, 980, 580, 55, 1000)
$row18 = @(2.46, 2.62, 3.25, 3.6, 3, 3.25, 1.54, 1.01, 2.68, 1.47, 1.58, 2.46, 1.19, 2.48, 1.84, 1.85, 1.38, 1.61, 9.6, 10.5, 42, 1000, 8.199999999999999, 7.4, 15.5, 230, 15.5, 12.5, 65, 1000, 220, 170, 420, 580, 1000, 1000)
$row19 = @(3, 3.25, 2.5, 2.68, 3.4, 3.45, 1.42, 1.07, 3.5, 1.34, 1.9, 2.02, 1.33, 3.55, 1.78, 2.16, 1.6, 1.45, 14, 10.5, 17, 980, 12, 7.8, 12.5, 980, 980, 14, 18, 980, 150, 38, 980, 580, 1000, 1000)
$row20 = @(3.95, 4.5, 2.02, 2.16, 3.4, 3.65, 1.01, 1.08, 2.8, 1.39, 1.74, 1.95, 1.27, 3.2, 1.92, 2, 1.86, 1.29, 25, 17.5, 980, 980, 90, 14, 40, 980, 980, 980, 980, 980, 1000, 300, 1000, 580, 1000, 55)
$row21 = @(2.48, 2.6, 3.25, 3.4, 3.2, 3.3, 1.01, 1.08, 2.88, 1.37, 1.84, 1.92, 1.29, 3.85, 1.82, 2.1, 1.42, 1.62, 22, 22, 70, 1000, 18.5, 7.6, 28, 1000, 1000, 23, 1000, 1000, 1000, 1000, 1000, 580, 1000, 1000)
$row22 = @(3.35, 3.7, 2.26, 2.5, 3.25, 3.5, 1.46, 1.08, 3.15, 1.39, 1.72, 1.95, 1.28, 3.45, 1.83, 1.96, 1.68, 1.37, 12.5, 9.4, 15.5, 40, 12.5, 8, 12, 32, 26, 15.5, 20, 48, 80, 250, 340, 580, 55, 26)
$row23 = @(1.71, 1.94, 4.3, 5.4, 3.55, 5.5, 1.01, 1.04, 4.4, 1.23, 2.2, 1.71, 1.47, 2.78, 1.68, 2.2, 1.22, 2.06, 1000, 1000, 1000, 1000, 1000, 42, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 85, 1000)
$row24 = @(2.14, 2.86, 2.7, 3.9, 3.45, 3.85, 1.01, 1.01, 4.1, 1.25, 1.25, 1.25, 1.32, 2.1, 1.03, 1.03, 1.34, 1.53, 1000, 1000, 1000, 1000, 1000, 42, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000)
$row25 = @(2.74, 2.76, 2.84, 2.88, 3.45, 3.5, 1.39, 1.06, 4.1, 1.3, 2.04, 1.93, 1.42, 3.3, 1.72, 2.3, 1.53, 1.56, 15.5, 12, 18.5, 42, 12, 7.6, 12.5, 30, 18, 12, 15.5, 40, 40, 28, 40, 80, 23, 24)
$row26 = @(1.66, 1.72, 6.2, 6.8, 3.7, 4.1, 1.48, 1.09, 3.05, 1.44, 1.68, 2.26, 1.24, 4.6, 2.18, 1.72, 1.17, 2.38, 11.5, 18.5, 55, 240, 6.8, 8.6, 28, 130, 8.800000000000001, 10.5, 29, 150, 17, 22, 55, 580, 15, 1000)
$row27 = @(2.9, 3, 2.52, 2.6, 3.6, 3.7, 1.37, 1.06, 4.2, 1.29, 2.1, 1.85, 1.43, 3.15, 1.7, 2.34, 1.62, 1.5, 17, 12, 17.5, 36, 13.5, 8.199999999999999, 12, 27, 21, 13.5, 16.5, 38, 48, 32, 40, 200, 26, 20)
$row28 = @(1.87, 2.06, 4.3, 5, 3.35, 3.85, 1.44, 1.01, 3.25, 1.37, 1.78, 2.06, 1.29, 3.75, 1.87, 1.94, 1.25, 1.94, 13, 15.5, 110, 130, 8.4, 8.199999999999999, 19.5, 1000, 12, 11, 21, 370, 29, 70, 120, 580, 17, 1000)

$rows = @($row2, $row3, $row4, $row5, $row6, $row7, $row8, $row9, $row10, $row11, $row12, $row13, $row14, $row15, $row16, $row17, $row18, $row19, $row20, $row21, $row22, $row23, $row24, $row25, $row26, $row27, $row28)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $rowIndex = $i + 2
    $rowData = $rows[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $colIndex = $j + 6
        $ws.Cells.Item($rowIndex, $colIndex).Value = $rowData[$j]
    }
}